$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2 through 109 in column C currently hold 7310 (Fitness values for
# early generations of run 6). Update them to 7293 to match the rest of
# the log (rows 110+ already report 7293).
for ($r = 2; $r -le 109; $r++) {
    $ws.Cells.Item($r, 3).Value = 7293
}
